# Add a new customer row ("Nguyễn Ngọc Hân") at the top of the data
# (row 2), pushing all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (first data row),
# shifting rows 2..51 down to 3..52.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new customer's data.
$ws.Range("A2").Value = "KH"
$ws.Range("B2").Value = 454
$ws.Range("C2").Value = "Nguyễn Ngọc Hân"
$ws.Range("D2").Value = "CẦN THƠ"
$ws.Range("F2").Value = "0587025000"
$ws.Range("I2").Value = 8500000
$ws.Range("J2").Value = 12000000
